$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.793.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.41%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.634.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.05%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.26%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'215.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.47%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.51%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.30%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2572"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.06407"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.74%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.36%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.61%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.285"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.55%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.862.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.05%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.630.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.64%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5602"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.59%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅7619"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.47%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'25.830.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.14%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'194.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.19%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.335"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.43%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.896"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.19%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.092"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.70%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.19%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.781"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -6.56%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'139.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.77%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.1253"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.19%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.803"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.14%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.77%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.242"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.09%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.04915"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.67%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.304"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.56%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.234"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.28%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.573"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.97%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.381"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.06%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.9026"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.37%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.572"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5533"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.72%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.126.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.29%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.01%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.03%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.498"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.31%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.8003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.49%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'98.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.53%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.772.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0₈111"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'55.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.21%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.4265"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.18%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.706"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.05%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05030"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.27%  "
$ws.Range("E50").Style = "Normal"
